$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")  # "Tracking Table" worksheet (internal name Sheet1)

# Insert a new row before row 5 (shifts existing rows 5-18 down to 6-19)
$ws.Rows.Item(5).Insert()

# Change F4 from 20 to 22 (top speed for Thornycroft J)
$ws.Range("F4").Value = 22

# Populate the new row 5: Leyland X Type
$ws.Range("A5").Value = "Leyland X Type"
$ws.Range("B5").Value = 1907
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "Heavy Goods"
$ws.Range("E5").Formula = "=IF(B5 > 1900, ((B5-1900)*10)+400+C5, ((B5-1730)*2)+C5)+VLOOKUP(D5,'ID Scheme'!`$A`$2:`$B`$4,2)"
$ws.Range("F5").Value = 20
$ws.Range("G5").Value = 18
$ws.Range("H5").Formula = "=SQRT(F5*G5)/`$B`$1"
$ws.Range("I5").Formula = "=H5*0.9"
$ws.Range("J5").Value = "x"

# Update the selected cell to F5, matching the captured view state
$ws.Range("F5").Select()
